# Generate Report for Handoff
#
# Updates the localization-status workbook to reflect a fresh handoff run for
# the "bb7ca973-2619-41a2-a4d2-ed0759aca5bf.md" source file (and the zh-cn
# status of "a10497b9-b8eb-41f8-bd43-444829d73014.md"):
#   - Status flips from "Handed back: in sync with en-US" to "Ready for handoff"
#   - The zh-cn row for bb7ca973 gets a new handoff timestamp and an error
#     noting the handback file is stale, with the Error Detail column widened
#     so the message is readable
#   - The de-de row for bb7ca973 gets the same treatment (new handoff
#     timestamp + stale-handback error + widened Error Detail column)
#   - The Overview sheet's summary row for bb7ca973 mirrors the new status
#     and the latest handoff-xliff-generate timestamp

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/b642be4af845834b076349f5b6e0b63fe428df40/e2e/bb7ca973-2619-41a2-a4d2-ed0759aca5bf.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/cd79cb5065d6f9ceb5d5edacab2cc15962b66280/e2e/bb7ca973-2619-41a2-a4d2-ed0759aca5bf.md."

# --- Overview sheet: bb7ca973 row (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $readyForHandoff
$overview.Range("F3").Value = $readyForHandoff
$overview.Range("G3").Value = "2016-08-12 16:59:51"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
# a10497b9 row (row 2): status only
$zhcn.Range("C2").Value = $readyForHandoff
# bb7ca973 row (row 3): status, new handoff datetime, error detail
$zhcn.Range("C3").Value = $readyForHandoff
$zhcn.Range("H3").Value = "2016-08-12 16:59:44"
$zhcn.Range("P3").Value = $errorDetail
# Widen the Error Detail column so the new message is legible (stored column
# width ends up 5/6 wider than the ColumnWidth we set, so back it off to land
# on an even 40)
$zhcn.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
# bb7ca973 row (row 3): status, new handoff datetime, error detail
$dede.Range("C3").Value = $readyForHandoff
$dede.Range("H3").Value = "2016-08-12 16:59:51"
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664
